{"js": "const segments = [\n  \"1. Defini\u00e7\u00e3o e classifica\u00e7\u00e3o de pol\u00edmeros termopl\u00e1sticos, elast\u00f4meros e fibras.\",\n  \"2. Identifica\u00e7\u00e3o de pl\u00e1sticos, borrachas e fibras.\",\n  \"3. Pol\u00edmeros de adi\u00e7\u00e3o olef\u00ednicos: polietileno, polipropileno e seus copol\u00edmeros.\",\n  \"4. Pol\u00edmeros de condensa\u00e7\u00e3o: poli(tereftalato de etileno), poli(tereftalato de butileno) e poliamidas.\",\n  \"5. Pol\u00edmeros halogenados: poli(cloreto de vinila), poli(tetrafluor etileno) e poli(fluoreto de vinilideno).\",\n  \"6. Termopl\u00e1sticos acr\u00edlicos e oximetil\u00eanicos: PMMA, POM e poliacetais.\",\n  \"7. Termopl\u00e1sticos nitrogenados: poliacrilonitrila, poliuretano, ABS e SAN.\",\n  \"8. Termopl\u00e1sticos estir\u00eanicos e fen\u00f3licos: poliestireno, HIPS, SBR e policarbonato.\",\n  \"9. Pol\u00edmeros hidrolis\u00e1veis: EVA, PVAc e PEO.\",\n  \"10. Termopl\u00e1sticos avan\u00e7ados: PPO, PPS e PEEK.\",\n  \"11. Elast\u00f4meros: borracha natural, polibutadieno, borrachas nitr\u00edlicas e fluoradas, EPDM e polisiloxanos.\",\n  \"12. Aditivos e compostos.\",\n  \"13. Tecnologias de transforma\u00e7\u00e3o apropriadas a cada tipo de pl\u00e1stico: extrus\u00e3o, inje\u00e7\u00e3o, lamina\u00e7\u00e3o, calandragem, termoforma\u00e7\u00e3o e moldagem por sopro.\",\n  \"14. Testes e ensaios de pol\u00edmeros termopl\u00e1sticos e elast\u00f4meros.\",\n  \"15. Reciclagem.\"\n];\nconst joined = segments.join(\"\\v\");\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"1. Defini\u00e7\u00e3o e classifica\u00e7\u00e3o de pol\u00edmeros termopl\u00e1sticos\") === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Target paragraph not found\");\n}\n\n// Replace the paragraph's whole text with the same text split into 15\n// numbered items separated by manual line breaks (w:br), matching the\n// vertical-tab-delimited convention already used elsewhere in this\n// document (e.g. the \"Cr\u00e9ditos-aula\" bullet list).\nconst range = target.getRange(Word.RangeLocation.whole);\nrange.insertText(joined, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The 15 numbered items of the \"Programa\" section, currently concatenated\n# into a single run of plain text with no separators between them.\n$segments = @(\n    \"1. Defini\u00e7\u00e3o e classifica\u00e7\u00e3o de pol\u00edmeros termopl\u00e1sticos, elast\u00f4meros e fibras.\",\n    \"2. Identifica\u00e7\u00e3o de pl\u00e1sticos, borrachas e fibras.\",\n    \"3. Pol\u00edmeros de adi\u00e7\u00e3o olef\u00ednicos: polietileno, polipropileno e seus copol\u00edmeros.\",\n    \"4. Pol\u00edmeros de condensa\u00e7\u00e3o: poli(tereftalato de etileno), poli(tereftalato de butileno) e poliamidas.\",\n    \"5. Pol\u00edmeros halogenados: poli(cloreto de vinila), poli(tetrafluor etileno) e poli(fluoreto de vinilideno).\",\n    \"6. Termopl\u00e1sticos acr\u00edlicos e oximetil\u00eanicos: PMMA, POM e poliacetais.\",\n    \"7. Termopl\u00e1sticos nitrogenados: poliacrilonitrila, poliuretano, ABS e SAN.\",\n    \"8. Termopl\u00e1sticos estir\u00eanicos e fen\u00f3licos: poliestireno, HIPS, SBR e policarbonato.\",\n    \"9. Pol\u00edmeros hidrolis\u00e1veis: EVA, PVAc e PEO.\",\n    \"10. Termopl\u00e1sticos avan\u00e7ados: PPO, PPS e PEEK.\",\n    \"11. Elast\u00f4meros: borracha natural, polibutadieno, borrachas nitr\u00edlicas e fluoradas, EPDM e polisiloxanos.\",\n    \"12. Aditivos e compostos.\",\n    \"13. Tecnologias de transforma\u00e7\u00e3o apropriadas a cada tipo de pl\u00e1stico: extrus\u00e3o, inje\u00e7\u00e3o, lamina\u00e7\u00e3o, calandragem, termoforma\u00e7\u00e3o e moldagem por sopro.\",\n    \"14. Testes e ensaios de pol\u00edmeros termopl\u00e1sticos e elast\u00f4meros.\",\n    \"15. Reciclagem.\"\n)\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"1. Defini\")) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Target paragraph not found\"\n}\n\n# Re-join the items with manual line breaks (vertical tab == Word's\n# \"\\v\" / ^l manual-line-break char) so they land as <w:br/> elements\n# inside the run instead of separate paragraphs.\n$target.Range.Text = ($segments -join \"`v\")\n\n"}
